# Align "Third Party C Compiler" box with the other compiler boxes, and
# update the date fields (datetime1 field cache) on the handout/notes masters.
#
# All geometry in the underlying OOXML is stored in EMUs; the PowerPoint
# object model works in points (1 pt = 12700 EMU), so convert accordingly.

$EMU_PER_PT = 12700

function EmuToPt([double]$emu) {
    return $emu / $EMU_PER_PT
}

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Move the "Third Party C Compiler" rounded rectangle (shape id 116)
#    up so it lines up with the other "Third Party ..." compiler boxes.
# ---------------------------------------------------------------------
$s = $p.Slides.Item(3)

$compiler = $null
$connStart = $null
$connEnd = $null

for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)
    if ($shp.Id -eq 116) { $compiler = $shp }
    if ($shp.Id -eq 45) { $connStart = $shp }
    if ($shp.Id -eq 47) { $connEnd = $shp }
}

# Move the compiler box up by 660,853 EMU (new top = 1,404,133 EMU).
$compiler.Left = EmuToPt 8313956
$compiler.Top = EmuToPt 1404133
$compiler.Width = EmuToPt 1599805
$compiler.Height = EmuToPt 385889

# Connector 45 runs from the bottom of the compiler box (stCxn) down to a
# free-floating arrowhead; update both ends to follow the new box position.
$connStart.Left = EmuToPt 9113859
$connStart.Top = EmuToPt 1790022
$connStart.Width = EmuToPt 11939
$connStart.Height = EmuToPt 1481488

# Connector 47 runs from the "C Code" box down to the top of the compiler
# box (endCxn); only its length changes since the top stays fixed.
$connEnd.Left = EmuToPt 9110094
$connEnd.Top = EmuToPt 964985
$connEnd.Width = EmuToPt 3765
$connEnd.Height = EmuToPt 439148

# ---------------------------------------------------------------------
# 2) Bump the cached date field text from 06/11/2025 to 07/11/2025 on the
#    handout master and the notes master.
# ---------------------------------------------------------------------
$handoutMaster = $p.HandoutMaster
for ($i = 1; $i -le $handoutMaster.Shapes.Count; $i++) {
    $shp = $handoutMaster.Shapes.Item($i)
    if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
        $tr = $shp.TextFrame.TextRange
        if ($tr.Text -eq "06/11/2025") {
            $tr.Text = "07/11/2025"
        }
    }
}

$notesMaster = $p.NotesMaster
for ($i = 1; $i -le $notesMaster.Shapes.Count; $i++) {
    $shp = $notesMaster.Shapes.Item($i)
    if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
        $tr = $shp.TextFrame.TextRange
        if ($tr.Text -eq "06/11/2025") {
            $tr.Text = "07/11/2025"
        }
    }
}
